$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on Price cells whose new values look numeric,
# so Excel does not silently coerce them into floating-point numbers.
$textCells = @("D5","D6","D11","D12","D13","D14","D18","D20","D22","D23","D24","D25","D27","D28","D30","D31","D32","D33","D35","D39","D40","D41","D42","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "71.204.86"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "3.839.67"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "711.93"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").Value = "173.12"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").Value = "3.838.93"
$ws.Range("E7").Value = "  +0.94%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").Value = "7.32"
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("D13").Value = "0.0000257"
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").Value = "36.82"
$ws.Range("E14").Value = "  +1.65%  "
$ws.Range("D15").Value = "4.486.26"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").Value = "3.812.96"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "71.149.91"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "7.23"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("D20").Value = "17.40"
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("E21").Value = "  -3.92%  "
$ws.Range("D22").Value = "495.72"
$ws.Range("E22").Value = "  +3.47%  "
$ws.Range("D23").Value = "0.725"
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("D24").Value = "85.16"
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("D25").Value = "0.0000147"
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("E26").Value = "  +1.98%  "
$ws.Range("D27").Value = "12.16"
$ws.Range("E27").Value = "  -1.44%  "
$ws.Range("D28").Value = "3.21"
$ws.Range("E28").Value = "  +2.72%  "
$ws.Range("E29").Value = "  -2.72%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").Value = "7.53"
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("D32").Value = "2.25"
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("D33").Value = "29.51"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E34").Value = "  -5.15%  "
$ws.Range("D35").Value = "9.22"
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("D36").Value = "3.804.07"
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "6.03"
$ws.Range("E39").Value = "  +0.48%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").Value = "1.03"
$ws.Range("E40").Value = "  +5.56%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "3.36"
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "2.28"
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "0.000316"
$ws.Range("E45").Value = "  -3.44%  "
$ws.Range("D46").Value = "163.65"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").Value = "48.85"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "417.52"
$ws.Range("E48").Value = "  +2.08%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "1.39"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("D50").Value = "8.64"
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("D51").Value = "0.297"
$ws.Range("E51").Value = "  -1.15%  "
